# them sap xep mang cau truc,having sql...
# Append 14 new rows (930-943) of vocabulary entries to Sheet1, columns A:T,
# mirroring the repeated-row pattern already used throughout this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A930:A943").Value = "accuracy chinh xac"
$ws.Range("B930:B943").Value = "acceptable co the chap nhan duoc"
$ws.Range("C930:C943").Value = "precise chinh xac"
$ws.Range("D930:D943").Value = "revelant thich hop co lien quan"
$ws.Range("E930:E943").Value = "inpiration su cam hung"
$ws.Range("F930:F943").Value = "intersection giao diem"
$ws.Range("G930:G943").Value = "respective tuong ung"
$ws.Range("H930:H943").Value = "schema luoc do"
$ws.Range("I930:I943").Value = "rigid cung"
$ws.Range("J930:J943").Value = "alternative su thay the"
$ws.Range("K930:K943").Value = "beam trum"
$ws.Range("L930:L943").Value = "inertia quan tinh"
$ws.Range("M930:M943").Value = "noticeable de nhan thay"
$ws.Range("N930:N943").Value = "prediction su tien doan"
$ws.Range("O930:O943").Value = "quantity so luong"
$ws.Range("P930:P943").Value = "establish thiet lap"
$ws.Range("Q930:Q943").Value = "permanent vinh vien"
$ws.Range("R930:R943").Value = "diverse nhieu loai"
$ws.Range("S930:S943").Value = "shape hinh dang"
$ws.Range("T930:T943").Value = "virtual ao"

# Move the viewport/selection to the new bottom-right corner of the sheet,
# matching where Excel would land after typing this block of rows.
$ws.Application.ActiveWindow.ScrollRow = 925
$ws.Application.ActiveWindow.ScrollColumn = 17
$ws.Range("U943").Select()
